# Cosmetic updates to DOCX template files:
#  1. Add the xmlns:oel="http://schemas.microsoft.com/office/2019/extlst"
#     namespace declaration (right after xmlns:o="...office:office") to the
#     root elements of document.xml, endnotes.xml, footer1-3.xml,
#     footnotes.xml and header1-3.xml.
#  2. Add a few new <w:lsdException/> entries inside the <w:latentStyles/>
#     block of styles.xml.
#
# We operate on the whole flat-OPC WordOpenXML blob so that every affected
# part (document/header/footer/footnotes/endnotes/styles) is touched in one
# pass; the search strings below are crafted to be unique within that blob
# (they deliberately do NOT match word/settings.xml, which the diff leaves
# untouched).

$d = $word.ActiveDocument
$xml = $d.WordOpenXML

# --- 1. namespace declaration -------------------------------------------
# This exact preamble fragment (long chartex/aink/am3d run immediately
# followed by xmlns:o then xmlns:r) occurs exactly once per affected part
# (document.xml, endnotes.xml, footer1.xml, footer2.xml, footer3.xml,
# footnotes.xml, header1.xml, header2.xml, header3.xml) -- 9 times total --
# and .Replace() rewrites every occurrence.
$oldNs = 'xmlns:am3d="http://schemas.microsoft.com/office/drawing/2017/model3d" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"'
$newNs = 'xmlns:am3d="http://schemas.microsoft.com/office/drawing/2017/model3d" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:oel="http://schemas.microsoft.com/office/2019/extlst" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"'
$xml = $xml.Replace($oldNs, $newNs)

# --- 2. styles.xml latent style exceptions -------------------------------

# a) "Normal Table" inserted right before "annotation subject"
$old1 = '<w:lsdException w:name="HTML Variable" w:semiHidden="1" w:unhideWhenUsed="1"/><w:lsdException w:name="annotation subject" w:semiHidden="1" w:unhideWhenUsed="1"/>'
$new1 = '<w:lsdException w:name="HTML Variable" w:semiHidden="1" w:unhideWhenUsed="1"/><w:lsdException w:name="Normal Table" w:semiHidden="1" w:unhideWhenUsed="1"/><w:lsdException w:name="annotation subject" w:semiHidden="1" w:unhideWhenUsed="1"/>'
$xml = $xml.Replace($old1, $new1)

# b) "Table Subtle 1" inserted right before "Table Subtle 2"
$old2 = '<w:lsdException w:name="Table Professional" w:semiHidden="1" w:unhideWhenUsed="1"/><w:lsdException w:name="Table Subtle 2" w:semiHidden="1" w:unhideWhenUsed="1"/>'
$new2 = '<w:lsdException w:name="Table Professional" w:semiHidden="1" w:unhideWhenUsed="1"/><w:lsdException w:name="Table Subtle 1" w:semiHidden="1" w:unhideWhenUsed="1"/><w:lsdException w:name="Table Subtle 2" w:semiHidden="1" w:unhideWhenUsed="1"/>'
$xml = $xml.Replace($old2, $new2)

# c) "Table Web 2" and "Table Web 3" inserted right after "Table Web 1"
$old3 = '<w:lsdException w:name="Table Web 1" w:semiHidden="1" w:unhideWhenUsed="1"/><w:lsdException w:name="Balloon Text" w:semiHidden="1" w:unhideWhenUsed="1"/>'
$new3 = '<w:lsdException w:name="Table Web 1" w:semiHidden="1" w:unhideWhenUsed="1"/><w:lsdException w:name="Table Web 2" w:semiHidden="1" w:unhideWhenUsed="1"/><w:lsdException w:name="Table Web 3" w:semiHidden="1" w:unhideWhenUsed="1"/><w:lsdException w:name="Balloon Text" w:semiHidden="1" w:unhideWhenUsed="1"/>'
$xml = $xml.Replace($old3, $new3)

$d.WordOpenXML = $xml

Write-Output "done"
